$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E1").Value = "Gini Rasio"

$ws.Range("E2").Value = 0.29399999999999998
$ws.Range("E3").Value = 0.30149999999999999
$ws.Range("E4").Value = 0.28499999999999998
$ws.Range("E5").Value = 0.30649999999999999
$ws.Range("E6").Value = 0.318
$ws.Range("E7").Value = 0.33200000000000002
$ws.Range("E8").Value = 0.34250000000000003
$ws.Range("E9").Value = 0.30149999999999999
$ws.Range("E10").Value = 0.23949999999999999
$ws.Range("E11").Value = 0.35299999999999998
$ws.Range("E12").Value = 0.42699999999999999
$ws.Range("E13").Value = 0.42449999999999999
$ws.Range("E14").Value = 0.36549999999999999
$ws.Range("E15").Value = 0.43149999999999999
$ws.Range("E16").Value = 0.3725
$ws.Range("E17").Value = 0.35599999999999998
$ws.Range("E18").Value = 0.35449999999999998
$ws.Range("E19").Value = 0.36249999999999999
$ws.Range("E20").Value = 0.316
$ws.Range("E21").Value = 0.312
$ws.Range("E22").Value = 0.30249999999999999
$ws.Range("E23").Value = 0.3
$ws.Range("E24").Value = 0.3155
$ws.Range("E25").Value = 0.26150000000000001
$ws.Range("E26").Value = 0.35349999999999998
$ws.Range("E27").Value = 0.30499999999999999
$ws.Range("E28").Value = 0.36149999999999999
$ws.Range("E29").Value = 0.36749999999999999
$ws.Range("E30").Value = 0.41349999999999998
$ws.Range("E31").Value = 0.34199999999999997
$ws.Range("E32").Value = 0.28649999999999998
$ws.Range("E33").Value = 0.30599999999999999
$ws.Range("E34").Value = 0.38700000000000001
$ws.Range("E35").Value = 0.34649999999999997
$ws.Range("E36").Value = 0.38350000000000001
$ws.Range("E37").Value = 0.41400000000000003
$ws.Range("E38").Value = 0.36799999999999999
$ws.Range("E39").Value = 0.34299999999999997

$ws.Range("G5").Select()

Write-Host "done"
